$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "statmech_model" column entries (F8, F11) from "Harmonic" to
# "Electronic" first so the new shared string "Electronic" is appended
# before the updated description text below (matches authoring order).
$ws.Range("F8").Value = "Electronic"
$ws.Range("F11").Value = "Electronic"

# Update the description cell for the statmech_model header (F2) to point
# to the new presets documentation instead of the old IdealGas/Harmonic text.
$ws.Range("F2").Value = "Type of thermodynamic model. See presets available here: https://vlachosgroup.github.io/pMuTT/statmech.html?highlight=presets#presets"
